$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Sub Category " header to "sub_category"
$ws.Range("B2").Value = "sub_category"

# Update the active cell / selection on the sheet
$ws.Range("C24").Select()
